$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly data rows are inserted above the existing block that starts
# at row 339, pushing all the existing rows (339..430) down by two rows
# (to 341..432). The two new rows hold a fresh "Especial"/"Primera" pair of
# observations (week of 2023-01-06, serial 44932) for Femacal de La Calera -
# Frutilla.
$ws.Rows("339:340").Insert()

# New row 339: Especial
$ws.Cells.Item(339, 1).Value = 3
$ws.Cells.Item(339, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(339, 3).Value = "Coquimbo"
$ws.Cells.Item(339, 4).Value = 44932
$ws.Cells.Item(339, 5).Value = 5
$ws.Cells.Item(339, 6).Value = "Fruta"
$ws.Cells.Item(339, 7).Value = 100101
$ws.Cells.Item(339, 8).Value = "Berries"
$ws.Cells.Item(339, 9).Value = 100112025
$ws.Cells.Item(339, 10).Value = "Frutilla"
$ws.Cells.Item(339, 11).Value = "Sin especificar"
$ws.Cells.Item(339, 12).Value = "Especial"
$ws.Cells.Item(339, 13).Value = 56
$ws.Cells.Item(339, 14).Value = 9000
$ws.Cells.Item(339, 15).Value = 9000
$ws.Cells.Item(339, 16).Value = 9000
$ws.Cells.Item(339, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(339, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(339, 19).Value = 1286
$ws.Cells.Item(339, 20).Value = 7

# New row 340: Primera
$ws.Cells.Item(340, 1).Value = 3
$ws.Cells.Item(340, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(340, 3).Value = "Coquimbo"
$ws.Cells.Item(340, 4).Value = 44932
$ws.Cells.Item(340, 5).Value = 5
$ws.Cells.Item(340, 6).Value = "Fruta"
$ws.Cells.Item(340, 7).Value = 100101
$ws.Cells.Item(340, 8).Value = "Berries"
$ws.Cells.Item(340, 9).Value = 100112025
$ws.Cells.Item(340, 10).Value = "Frutilla"
$ws.Cells.Item(340, 11).Value = "Sin especificar"
$ws.Cells.Item(340, 12).Value = "Primera"
$ws.Cells.Item(340, 13).Value = 48
$ws.Cells.Item(340, 14).Value = 7000
$ws.Cells.Item(340, 15).Value = 7000
$ws.Cells.Item(340, 16).Value = 7000
$ws.Cells.Item(340, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(340, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(340, 19).Value = 1000
$ws.Cells.Item(340, 20).Value = 7
